# Add new trading-journal entries (rows 127-132) to the "intraday" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intraday")

# Reuse the date format (style) that column A already uses for its date
# cells by copying it from the last existing data row down onto the new
# rows, then overwrite each cell's value.
$ws.Range("A126").Copy()
$ws.Range("A127:A132").PasteSpecial(-4122)

# Row 127 - 2025-03-13
$ws.Range("A127").Value = 45729
$ws.Range("B127").Value = " "
$ws.Range("C127").Value = 473
$ws.Range("D127").Value = 1
$ws.Range("E127").Value = "trade was according to my setup but it went up after hitting my SL "

# Row 128 - 2025-03-17
$ws.Range("A128").Value = 45733
$ws.Range("B128").Value = 509
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = "trade was according to my setup so went well 1:2.5 risk to reqard I captured "

# Row 129 - 2025-03-18
$ws.Range("A129").Value = 45734
$ws.Range("B129").Value = 23
$ws.Range("D129").Value = 1
$ws.Range("E129").Value = "was a good trade as per my setup but due to time issue I exited early and trade gave amlost 80 points "

# Row 130 - note-only row (no date); clear the pasted date style/format
# since this row has no date value in column A.
$ws.Range("A130").Clear()
$ws.Range("B130").Value = "If you took trade after 3pm wait till 3:25 that works in option "

# Row 131 - 2025-03-19
$ws.Range("A131").Value = 45735
$ws.Range("B131").Value = 0
$ws.Range("D131").Value = 0
$ws.Range("E131").Value = "No trade as I didn't find any move in market "

# Row 132 - 2025-03-20
$ws.Range("A132").Value = 45736
$ws.Range("B132").Value = 745
$ws.Range("D132").Value = 1
$ws.Range("E132").Value = "perfect trade as per gap up and  captured small profit in 1 min setup "

# Match the author's final on-screen selection/scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 102
$null = $ws.Range("D132").Select()
